# Fixed sample products cupcake typo
#
# The "Cup Cake Sample" product's name/meta text had an inconsistent
# "Cup Cake" spelling; this normalizes it to "Cupcake" to match the
# slug/SKU naming used elsewhere, and updates the Categories sheet's
# reference slug to the corrected value.

$wb = $excel.ActiveWorkbook

# --- Main sheet: row 6 is the Cup Cake Sample product -------------------
$main = $wb.Worksheets.Item("Main")

# Name column (E)
$main.Range("E6").Value = "Cupcake Sample"
# Meta Title column (O)
$main.Range("O6").Value = "Vanilla Cupcake with Rich Frosting"
# Meta Description column (P)
$main.Range("P6").Value = "Vanilla Cupcake with Rich Frosting"

# --- Categories sheet: update the matching slug reference ----------------
$cats = $wb.Worksheets.Item("Categories")
$cats.Range("A5").Value = "cupcake-sample"

# --- Make Categories the active/selected sheet+cell, matching the saved
#     workbook view state from the edit session ---------------------------
$cats.Activate() | Out-Null
$cats.Range("A5").Select() | Out-Null
